# Apply edits to the "axes" worksheet: insert three new columns (D, E, F)
# holding the arrow-axis labels / percentage labels, pushing the old
# "Title" / "Feldspar Classification Diagram" column out to column G.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("axes")

# Shift the existing 4th column (Title / Feldspar Classification Diagram)
# to column G before writing the new columns in D:F.
$ws.Range("D1").Value = "Title"
$ws.Range("D2").Value = "Feldspar Classification Diagram"
$ws.Range("D1:D2").Cut($ws.Range("G1"))

# New columns with the axis-arrow labels and percentage labels.
$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"

$ws.Range("D2").Value = "Or (%)"
$ws.Range("E2").Value = "Ab (%)"
$ws.Range("F2").Value = "An (%)"

# Widen the new columns to fit their content (target display width ~14.1
# characters).
$ws.Columns("D:F").ColumnWidth = 13.25

# Make "axes" the active sheet / tab, with F2 selected (matching the
# last cell touched while filling in the new columns).
$ws.Activate() | Out-Null
$ws.Range("F2").Select() | Out-Null
